$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '34.025.84'
$cell.Style = 'Normal'

$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  -0.16%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.776.27'
$cell.Style = 'Normal'

$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  -2.29%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '225.35'
$cell.Style = 'Normal'

$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  -1.28%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'

$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  +0.00%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '31.63'
$cell.Style = 'Normal'

$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  +0.80%  '
$cell.Style = 'Normal'

$ws.Range('B9').Value = 'Cardano'

$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.280'
$cell.Style = 'Normal'

$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  -0.32%  '
$cell.Style = 'Normal'

$ws.Range('B10').Value = 'Dogecoin'

$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.0656'
$cell.Style = 'Normal'

$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  -1.68%  '
$cell.Style = 'Normal'

$ws.Range('B11').Value = 'TRON'

$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'

$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.0928'
$cell.Style = 'Normal'

$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.Style = 'Normal'

$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'

$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'

$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '2.030.11'
$cell.Style = 'Normal'

$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  -2.33%  '
$cell.Style = 'Normal'

$ws.Range('B13').Value = 'Chainlink'

$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '10.96'
$cell.Style = 'Normal'

$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  +5.90%  '
$cell.Style = 'Normal'

$ws.Range('B14').Value = 'WrappedEther'

$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '1.769.23'
$cell.Style = 'Normal'

$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  -2.82%  '
$cell.Style = 'Normal'

$ws.Range('B15').Value = 'Polygon'

$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '0.623'
$cell.Style = 'Normal'

$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  -2.95%  '
$cell.Style = 'Normal'

$ws.Range('B16').Value = 'WrappedBTC'

$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '33.995.64'
$cell.Style = 'Normal'

$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  -0.26%  '
$cell.Style = 'Normal'

$ws.Range('B17').Value = 'Polkadot'

$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'

$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '4.20'
$cell.Style = 'Normal'

$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  -1.63%  '
$cell.Style = 'Normal'

$ws.Range('B18').Value = 'Litecoin'

$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '68.73'
$cell.Style = 'Normal'

$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  -1.05%  '
$cell.Style = 'Normal'

$ws.Range('B19').Value = 'BitcoinCash'

$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'

$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '252.89'
$cell.Style = 'Normal'

$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  -2.04%  '
$cell.Style = 'Normal'

$ws.Range('B20').Value = 'ShibaInu'

$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0738'
$cell.Style = 'Normal'

$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  -1.47%  '
$cell.Style = 'Normal'

$ws.Range('B21').Value = 'Dai'

$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'

$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  +0.19%  '
$cell.Style = 'Normal'

$ws.Range('B22').Value = 'Avalanche'

$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '10.35'
$cell.Style = 'Normal'

$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  -1.57%  '
$cell.Style = 'Normal'

$ws.Range('B23').Value = 'Uniswap'

$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '4.19'
$cell.Style = 'Normal'

$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  -3.17%  '
$cell.Style = 'Normal'

$ws.Range('B24').Value = 'Toncoin'

$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '2.14'
$cell.Style = 'Normal'

$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -2.70%  '
$cell.Style = 'Normal'

$ws.Range('B25').Value = 'Monero'

$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '156.43'
$cell.Style = 'Normal'

$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  -0.96%  '
$cell.Style = 'Normal'

$ws.Range('B26').Value = 'EthereumClassic'

$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '16.39'
$cell.Style = 'Normal'

$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  -0.99%  '
$cell.Style = 'Normal'

$ws.Range('B27').Value = 'Cosmos'

$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '6.98'
$cell.Style = 'Normal'

$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  -2.28%  '
$cell.Style = 'Normal'

$ws.Range('B28').Value = 'Stellar'

$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '0.113'
$cell.Style = 'Normal'

$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  -0.92%  '
$cell.Style = 'Normal'

$ws.Range('B29').Value = 'BinanceUSD'

$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'

$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'

$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  +0.00%  '
$cell.Style = 'Normal'

$ws.Range('B30').Value = 'Filecoin'

$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '3.76'
$cell.Style = 'Normal'

$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  -2.64%  '
$cell.Style = 'Normal'

$ws.Range('B31').Value = 'Hedera'

$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.0511'
$cell.Style = 'Normal'

$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  -0.39%  '
$cell.Style = 'Normal'

$ws.Range('B32').Value = 'PancakeSwap'

$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '1.19'
$cell.Style = 'Normal'

$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -0.68%  '
$cell.Style = 'Normal'

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'

$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '3.57'
$cell.Style = 'Normal'

$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  +1.65%  '
$cell.Style = 'Normal'

$ws.Range('B34').Value = 'LidoDAOToken'

$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '1.85'
$cell.Style = 'Normal'

$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  +2.54%  '
$cell.Style = 'Normal'

$ws.Range('B35').Value = 'Maker'

$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.445.86'
$cell.Style = 'Normal'

$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  -6.24%  '
$cell.Style = 'Normal'

$ws.Range('B36').Value = 'TrustWalletToken'

$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '1.05'
$cell.Style = 'Normal'

$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -3.02%  '
$cell.Style = 'Normal'

$ws.Range('B37').Value = 'VeChain'

$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.0187'
$cell.Style = 'Normal'

$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  -0.31%  '
$cell.Style = 'Normal'

$ws.Range('B38').Value = 'ImmutableX'

$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.623'
$cell.Style = 'Normal'

$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  +0.05%  '
$cell.Style = 'Normal'

$ws.Range('B39').Value = 'MXToken'

$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '2.84'
$cell.Style = 'Normal'

$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  +0.91%  '
$cell.Style = 'Normal'

$ws.Range('B40').Value = 'Aave'

$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'

$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '82.77'
$cell.Style = 'Normal'

$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  -2.54%  '
$cell.Style = 'Normal'

$ws.Range('B41').Value = 'HuobiToken'

$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '2.34'
$cell.Style = 'Normal'

$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  +0.11%  '
$cell.Style = 'Normal'

$ws.Range('B42').Value = 'ARBITRUM'

$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.890'
$cell.Style = 'Normal'

$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  -2.70%  '
$cell.Style = 'Normal'

$ws.Range('B43').Value = 'RenderToken'

$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '2.04'
$cell.Style = 'Normal'

$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  -5.10%  '
$cell.Style = 'Normal'

$ws.Range('B44').Value = 'Kaspa'

$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'

$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.0506'
$cell.Style = 'Normal'

$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  -2.73%  '
$cell.Style = 'Normal'

$ws.Range('B45').Value = 'WEMIXToken'

$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '1.05'
$cell.Style = 'Normal'

$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  -2.04%  '
$cell.Style = 'Normal'

$ws.Range('B46').Value = 'RocketPoolETH'

$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '1.931.35'
$cell.Style = 'Normal'

$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  -2.27%  '
$cell.Style = 'Normal'

$ws.Range('B47').Value = 'FraxShare'

$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '5.78'
$cell.Style = 'Normal'

$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  +1.07%  '
$cell.Style = 'Normal'

$ws.Range('B48').Value = 'InjectiveProtocol'

$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '12.05'
$cell.Style = 'Normal'

$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +3.38%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  +0.13%  '
$cell.Style = 'Normal'

$ws.Range('B50').Value = 'Quant'

$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '97.90'
$cell.Style = 'Normal'

$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  +2.83%  '
$cell.Style = 'Normal'

$ws.Range('B51').Value = 'BitcoinSV'

$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'

$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '49.57'
$cell.Style = 'Normal'

$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  -5.89%  '
$cell.Style = 'Normal'
